# Edit the sl fmri experiment stimuli (visual_run1_5.xlsx)
# Updates the trial table on Sheet1 (rows 2-49) with the new stimulus
# ordering/numbering, extends the table with 12 additional trial rows,
# and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{R=2;  A="M.png"; B=97;  C="S"; D=1; E=1; F=10},
    @{R=3;  A="K.png"; B=98;  C="S"; D=1; E=2; F=10},
    @{R=4;  A="B.png"; B=99;  C="S"; D=1; E=3; F=10},
    @{R=5;  A="A.png"; B=100; C="S"; D=2; E=1; F=8},
    @{R=6;  A="J.png"; B=101; C="S"; D=2; E=2; F=8},
    @{R=7;  A="F.png"; B=102; C="S"; D=2; E=3; F=8},
    @{R=8;  A="E.png"; B=103; C="S"; D=4; E=1; F=9},
    @{R=9;  A="L.png"; B=104; C="S"; D=4; E=2; F=9},
    @{R=10; A="H.png"; B=105; C="S"; D=4; E=3; F=9},
    @{R=11; A="D.png"; B=106; C="S"; D=3; E=1; F=9},
    @{R=12; A="G.png"; B=107; C="S"; D=3; E=2; F=9},
    @{R=13; A="C.png"; B=108; C="S"; D=3; E=3; F=9},
    @{R=14; A="A.png"; B=109; C="S"; D=2; E=1; F=9},
    @{R=15; A="J.png"; B=110; C="S"; D=2; E=2; F=9},
    @{R=16; A="F.png"; B=111; C="S"; D=2; E=3; F=9},
    @{R=17; A="M.png"; B=112; C="S"; D=1; E=1; F=11},
    @{R=18; A="K.png"; B=113; C="S"; D=1; E=2; F=11},
    @{R=19; A="B.png"; B=114; C="S"; D=1; E=3; F=11},
    @{R=20; A="D.png"; B=115; C="S"; D=3; E=1; F=10},
    @{R=21; A="G.png"; B=116; C="S"; D=3; E=2; F=10},
    @{R=22; A="C.png"; B=117; C="S"; D=3; E=3; F=10},
    @{R=23; A="E.png"; B=118; C="S"; D=4; E=1; F=10},
    @{R=24; A="L.png"; B=119; C="S"; D=4; E=2; F=10},
    @{R=25; A="H.png"; B=120; C="S"; D=4; E=3; F=10},
    @{R=26; A="A.png"; B=121; C="S"; D=2; E=1; F=10},
    @{R=27; A="J.png"; B=122; C="S"; D=2; E=2; F=10},
    @{R=28; A="F.png"; B=123; C="S"; D=2; E=3; F=10},
    @{R=29; A="M.png"; B=124; C="S"; D=1; E=1; F=12},
    @{R=30; A="K.png"; B=125; C="S"; D=1; E=2; F=12},
    @{R=31; A="B.png"; B=126; C="S"; D=1; E=3; F=12},
    @{R=32; A="A.png"; B=127; C="S"; D=2; E=1; F=11},
    @{R=33; A="J.png"; B=128; C="S"; D=2; E=2; F=11},
    @{R=34; A="F.png"; B=129; C="S"; D=2; E=3; F=11},
    @{R=35; A="E.png"; B=130; C="S"; D=4; E=1; F=11},
    @{R=36; A="L.png"; B=131; C="S"; D=4; E=2; F=11},
    @{R=37; A="H.png"; B=132; C="S"; D=4; E=3; F=11},
    @{R=38; A="D.png"; B=133; C="S"; D=3; E=1; F=11},
    @{R=39; A="G.png"; B=134; C="S"; D=3; E=2; F=11},
    @{R=40; A="C.png"; B=135; C="S"; D=3; E=3; F=11},
    @{R=41; A="E.png"; B=136; C="S"; D=4; E=1; F=12},
    @{R=42; A="L.png"; B=137; C="S"; D=4; E=2; F=12},
    @{R=43; A="H.png"; B=138; C="S"; D=4; E=3; F=12},
    @{R=44; A="D.png"; B=139; C="S"; D=3; E=1; F=12},
    @{R=45; A="G.png"; B=140; C="S"; D=3; E=2; F=12},
    @{R=46; A="C.png"; B=141; C="S"; D=3; E=3; F=12},
    @{R=47; A="A.png"; B=142; C="S"; D=2; E=1; F=12},
    @{R=48; A="J.png"; B=143; C="S"; D=2; E=2; F=12},
    @{R=49; A="F.png"; B=144; C="S"; D=2; E=3; F=12}
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
}

# Update the window scroll position / active selection to match the
# author's final view of the sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("C52").Select()
